$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(3.0,1,19.31916866666667,57.957506,0.09973928790435696,0.1012662650824037,3.0,1,10.28467933333333,30.854038,0.113655682908653,0.11797514492277,198.6914547232476,1788.223092509228,0.01133593687959245,0.01194690229888421),
    @(3.0,1,19.31916866666667,57.957506,0.09973928790435696,0.1012662650824037,3.0,1,1.769189333333333,5.307568,0.01955125827044465,0.02029430001957787,34.17926713393422,307.613404205408,0.00195002857752832,0.002055127965444402),
    @(3.0,1,19.31916866666667,57.957506,0.09973928790435696,0.1012662650824037,3.0,1,34.039953,102.119859,0.3761744998557516,0.3904709381967389,657.623593412406,5918.612340711654,0.0375193767433903,0.03954153353440581),
    @(3.0,1,19.31916866666667,57.957506,0.09973928790435696,0.1012662650824037,3.0,1,34.45657633333334,103.369729,0.3807785918192425,0.395250007775399,665.6724098595416,5991.051688735874,0.03797858559727504,0.04002549206120567),
    @(3.0,1,19.31916866666667,57.957506,0.09973928790435696,0.1012662650824037,2.0,1,9.939396,19.878792,0.1098399671459082,0.0760096090855142,192.020867768792,1152.125206612752,0.01095536010657085,0.007697209222463559),
    @(3.0,1,127.3682276666667,382.104683,0.6575653719009243,0.6676324912584373,3.0,1,10.28467933333333,30.854038,0.113655682908653,0.11797514492277,1309.941378806662,11789.47240925996,0.07473604140048197,0.07876403991136409),
    @(3.0,1,127.3682276666667,382.104683,0.6575653719009243,0.6676324912584373,3.0,1,1.769189333333333,5.307568,0.01955125827044465,0.02029430001957787,225.3385097934382,2028.046588140944,0.01285623041573596,0.01354913408041692),
    @(3.0,1,127.3682276666667,382.104683,0.6575653719009243,0.6676324912584373,3.0,1,34.039953,102.119859,0.3761744998557516,0.3904709381967389,4335.608483466633,39020.47635119969,0.2473593248972915,0.2606910852323081),
    @(3.0,1,127.3682276666667,382.104683,0.6575653719009243,0.6676324912584373,3.0,1,34.45657633333334,103.369729,0.3807785918192425,0.395250007775399,4388.673059037879,39498.05753134091,0.2503868163415304,0.2638817473610063),
    @(3.0,1,127.3682276666667,382.104683,0.6575653719009243,0.6676324912584373,2.0,1,9.939396,19.878792,0.1098399671459082,0.0760096090855142,1265.963252597156,7595.779515582937,0.07222695884588443,0.05074648467334179),
    @(3.0,1,18.657769,55.97330699999999,0.09632467245626405,0.09779937294404559,3.0,1,10.28467933333333,30.854038,0.113655682908653,0.11797514492277,191.8891712404073,1727.002541163666,0.01094784642896901,0.01153789519642981),
    @(3.0,1,18.657769,55.97330699999999,0.09632467245626405,0.09779937294404559,3.0,1,1.769189333333333,5.307568,0.01955125827044465,0.02029430001957787,33.00912589859733,297.082133087376,0.001883268549008405,0.001984769816253048),
    @(3.0,1,18.657769,55.97330699999999,0.09632467245626405,0.09779937294404559,3.0,1,34.039953,102.119859,0.3761744998557516,0.3904709381967389,635.1095798448569,5715.986218603712,0.03623488548500423,0.03818781290851424),
    @(3.0,1,18.657769,55.97330699999999,0.09632467245626405,0.09779937294404559,3.0,1,34.45657633333334,103.369729,0.3807785918192425,0.395250007775399,642.8828417582004,5785.945575823803,0.036678373135346,0.03865520291656317),
    @(3.0,1,18.657769,55.97330699999999,0.09632467245626405,0.09779937294404559,2.0,1,9.939396,19.878792,0.1098399671459082,0.0760096090855142,185.446954567524,1112.681727405144,0.01058029885793641,0.007433692106285319),
    @(3.0,1,19.589352,58.768056,0.1011341664177781,0.1026824987478506,3.0,1,10.28467933333333,30.854038,0.113655682908653,0.11797514492277,201.470203667792,1813.231833010128,0.01149447274960994,0.01211398267080982),
    @(3.0,1,19.589352,58.768056,0.1011341664177781,0.1026824987478506,3.0,1,1.769189333333333,5.307568,0.01955125827044465,0.02029430001957787,34.657272605312,311.915453447808,0.00197730020760011,0.002083869436348809),
    @(3.0,1,19.589352,58.768056,0.1011341664177781,0.1026824987478506,3.0,1,34.039953,102.119859,0.3761744998557516,0.3904709381967389,666.820621380456,6001.385592424103,0.03804409447053603,0.04009453162245869),
    @(3.0,1,19.589352,58.768056,0.1011341664177781,0.1026824987478506,3.0,1,34.45657633333334,103.369729,0.3807785918192425,0.395250007775399,674.9820025085362,6074.838022576824,0.03850972547337447,0.04058525842848535),
    @(3.0,1,19.589352,58.768056,0.1011341664177781,0.1026824987478506,2.0,1,9.939396,19.878792,0.1098399671459082,0.0760096090855142,194.706326911392,1168.237961468352,0.01110857351665756,0.007804856589747925),
    @(2.0,1,8.76216,17.52432,0.04523650132067659,0.03061937196726285,3.0,1,10.28467933333333,30.854038,0.113655682908653,0.11797514492277,90.11600586736002,540.69603520416,0.005141385449999683,0.003612324845282035),
    @(2.0,1,8.76216,17.52432,0.04523650132067659,0.03061937196726285,3.0,1,1.769189333333333,5.307568,0.01955125827044465,0.02029430001957787,15.50192000896,93.01152005376,0.0008844305205718587,0.0006213987211146844),
    @(2.0,1,8.76216,17.52432,0.04523650132067659,0.03061937196726285,3.0,1,34.039953,102.119859,0.3761744998557516,0.3904709381967389,298.2635145784799,1789.58108747088,0.01701681825952956,0.01195597489905205),
    @(2.0,1,8.76216,17.52432,0.04523650132067659,0.03061937196726285,3.0,1,34.45657633333334,103.369729,0.3807785918192425,0.395250007775399,301.91403488488,1811.48420930928,0.01722509127171653,0.01210230700813848),
    @(2.0,1,8.76216,17.52432,0.04523650132067659,0.03061937196726285,2.0,1,9.939396,19.878792,0.1098399671459082,0.0760096090855142,87.09057805536,348.36231222144,0.004968775818858948,0.002327366493675601)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowvals = $data[$r]
    for ($c = 0; $c -lt $rowvals.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 5).Value = $rowvals[$c]
    }
}
